$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARCHITECTURE")

# --- Fill in the "Hs" column (E) for building-use rows 2-19 ---
# Numeric entries
$ws.Range("E2").Value = 4
$ws.Range("E3").Value = 4
$ws.Range("E5").Value = 4
$ws.Range("E6").Value = 4
$ws.Range("E9").Value = 4
$ws.Range("E10").Value = 5
$ws.Range("E11").Value = 0
$ws.Range("E12").Value = 5
$ws.Range("E13").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("E17").Value = 5
$ws.Range("E18").Value = 5
$ws.Range("E19").Value = 5

# --- Add the new legend rows referencing the new sources (order matters
# for shared-string allocation order) ---
$ws.Range("B25").Value = "http://www.sciencedirect.com/science/article/pii/S037877881630442X"
$ws.Range("B25").HorizontalAlignment = -4131
$ws.Range("A25").Value = "4"

# Text ("4,0") entries for the Hs column - reference source 4
$ws.Range("E4").Value = "4,0"
$ws.Range("E7").Value = "4,0"
$ws.Range("E8").Value = "4,0"

$ws.Range("A26").Value = "5"
$ws.Range("B26").Value = "Average NTU"
$ws.Range("B26").HorizontalAlignment = -4131

# --- Restore selection similar to the authored workbook ---
$ws.Range("E12").Select()
